$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Complete row 9 (existing row) with the Bollinger / PriceChange / UpDown values ---
$ws.Cells.Item(9, 23).Value = 0                           # W9  Bollinger
$ws.Cells.Item(9, 24).Value = -0.21000099999999833        # X9  PriceChange
$ws.Cells.Item(9, 25).Value = "Down"                       # Y9  UpDown

# --- Append new row 10 with the full set of values ---
$ws.Cells.Item(10, 1).Value = 42653.87945601852            # A10 Date
$ws.Cells.Item(10, 2).Value = 13                            # B10 ScoreFinal
$ws.Cells.Item(10, 3).Value = "Buy"                         # C10 Verdict
$ws.Cells.Item(10, 4).Value = 32                            # D10 totalSentiment
$ws.Cells.Item(10, 5).Value = 6037                          # E10 wordCount
$ws.Cells.Item(10, 6).Value = 865                           # F10 sentenceCount
$ws.Cells.Item(10, 7).Value = 64                            # G10 posWordPercentage
$ws.Cells.Item(10, 8).Value = 34                            # H10 negWordPercentage
$ws.Cells.Item(10, 9).Value = 88                            # I10 posPhrasePercentage
$ws.Cells.Item(10, 10).Value = 11                           # J10 negPhrasePercentage
$ws.Cells.Item(10, 11).Value = 8894                         # K10 ElapsedMs
$ws.Cells.Item(10, 12).Value = 151                          # L10 posWordCount
$ws.Cells.Item(10, 13).Value = 80                           # M10 negWordCount
$ws.Cells.Item(10, 14).Value = 46                           # N10 positivePhraseCount
$ws.Cells.Item(10, 15).Value = 6                            # O10 negativePhraseCount
$ws.Cells.Item(10, 16).Value = "Noun"                       # P10 Method
$ws.Cells.Item(10, 17).Value = 47.96375473473072            # Q10 RSI
$ws.Cells.Item(10, 18).Value = 0.49                         # R10 PEG
$ws.Cells.Item(10, 19).Value = 0.0521                       # S10 200Moving%
$ws.Cells.Item(10, 19).NumberFormat = "0.00%"
$ws.Cells.Item(10, 20).Value = -0.0214                      # T10 50Moving%
$ws.Cells.Item(10, 20).NumberFormat = "0.00%"
$ws.Cells.Item(10, 21).Value = 2.25                         # U10 PriceBook
$ws.Cells.Item(10, 22).Value = "N/A"                        # V10 Dividend
$ws.Cells.Item(10, 23).Value = 0                            # W10 Bollinger
